# Apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force Excel to store the value as TEXT (not auto-converted to a
    # number/date) by entering it with a leading apostrophe, the same
    # way a user would via the keyboard. Resetting the Style back to
    # "Normal" afterwards clears the quote-prefix formatting flag that
    # Excel attaches to the cell, so the cell style is left untouched.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

$ws.Range('D2').Value = '29.306.06'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '1.873.27'
$ws.Range('E3').Value = '  +0.49%  '
Set-TextCell 'D4' '1.000'
Set-TextCell 'D5' '0.7109'
$ws.Range('E5').Value = '  +0.08%  '
Set-TextCell 'D6' '241.56'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  +0.08%  '
Set-TextCell 'D8' '0.3109'
$ws.Range('E8').Value = '  +0.34%  '
Set-TextCell 'D9' '0.07774'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E10').Value = '  +1.99%  '
Set-TextCell 'D11' '0.08383'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '1.865.44'
$ws.Range('E12').Value = '  -0.19%  '
Set-TextCell 'D13' '5.240'
$ws.Range('E13').Value = '  +0.38%  '
Set-TextCell 'D14' '0.7108'
$ws.Range('E14').Value = '  +0.38%  '
Set-TextCell 'D15' '91.13'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '29.326.03'
$ws.Range('E16').Value = '  +0.51%  '
Set-TextCell 'D17' '6.052'
$ws.Range('E17').Value = '  +2.28%  '
Set-TextCell 'D18' '0.000008195'
$ws.Range('E18').Value = '  +4.99%  '
Set-TextCell 'D19' '240.03'
$ws.Range('E19').Value = '  -1.28%  '
Set-TextCell 'D20' '13.19'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').Value = '2.118.70'
$ws.Range('E21').Value = '  +0.24%  '
Set-TextCell 'D22' '0.9997'
$ws.Range('E22').Value = '  -0.01%  '
Set-TextCell 'D23' '7.744'
$ws.Range('E23').Value = '  -1.46%  '
Set-TextCell 'D24' '1.001'
$ws.Range('E24').Value = '  +0.09%  '
Set-TextCell 'D25' '0.1584'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  -0.40%  '
Set-TextCell 'D27' '9.015'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('E28').Value = '  +0.21%  '
Set-TextCell 'D29' '1.510'
$ws.Range('E29').Value = '  +0.65%  '
Set-TextCell 'D30' '4.400'
$ws.Range('E30').Value = '  +0.27%  '
Set-TextCell 'D31' '1.291'
$ws.Range('E31').Value = '  -2.39%  '
Set-TextCell 'D32' '4.311'
$ws.Range('E32').Value = '  +1.42%  '
Set-TextCell 'D33' '0.05294'
$ws.Range('E33').Value = '  +2.68%  '
Set-TextCell 'D34' '1.940'
$ws.Range('E34').Value = '  +1.62%  '
Set-TextCell 'D35' '1.178'
$ws.Range('E35').Value = '  +1.28%  '
Set-TextCell 'D36' '0.7429'
$ws.Range('E36').Value = '  -6.42%  '
Set-TextCell 'D37' '2.707'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').Value = '1.227.50'
$ws.Range('E39').Value = '  +5.28%  '
Set-TextCell 'D40' '2.730'
$ws.Range('E40').Value = '  +1.02%  '
Set-TextCell 'D41' '6.549'
$ws.Range('E41').Value = '  +4.95%  '
Set-TextCell 'D42' '0.8847'
$ws.Range('E42').Value = '  -0.60%  '
Set-TextCell 'D43' '109.21'
$ws.Range('E43').Value = '  +6.39%  '
Set-TextCell 'D44' '72.44'
$ws.Range('E44').Value = '  -0.57%  '
Set-TextCell 'D45' '1.001'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '2.016.94'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D47' '0.5194'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D48' '1.794'
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('E49').Value = '  +2.51%  '
Set-TextCell 'D50' '9.380'
$ws.Range('E50').Value = '  +0.82%  '
Set-TextCell 'D51' '0.4307'
$ws.Range('E51').Value = '  +0.73%  '
